# Auto-generated: update crypto price/volume snapshot cells (D/E columns, rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.766.35"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "3.257.52"
$ws.Range("E3").Value = "  +2.44%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'604.98"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("E6").Value = "  +1.95%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.258.65"
$ws.Range("E8").Value = "  +2.50%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("D11").Value = "'5.91"
$ws.Range("E11").Value = "  +3.50%  "
$ws.Range("D12").Value = "'0.506"
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("D13").Value = "'0.0000272"
$ws.Range("E13").Value = "  +2.60%  "
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").Value = "3.796.07"
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("D16").Value = "66.790.68"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "3.264.39"
$ws.Range("E18").Value = "  +2.82%  "
$ws.Range("D19").Value = "'0.113"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("D20").Value = "'508.07"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("D22").Value = "'0.754"
$ws.Range("E22").Value = "  +2.92%  "
$ws.Range("D23").Value = "'8.08"
$ws.Range("E23").Value = "  -1.13%  "
$ws.Range("D24").Value = "'14.83"
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("D25").Value = "'86.29"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("D26").Value = "'0.166"
$ws.Range("E26").Value = "  +85.76%  "
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("D29").Value = "'9.12"
$ws.Range("E29").Value = "  -0.62%  "
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("D31").Value = "'6.92"
$ws.Range("E31").Value = "  -1.97%  "
$ws.Range("E32").Value = "  -7.41%  "
$ws.Range("D33").Value = "'28.27"
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("E35").Value = "  -4.12%  "
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("D37").Value = "0.0₃0812"
$ws.Range("E37").Value = "  +19.60%  "
$ws.Range("E38").Value = "  +18.96%  "
$ws.Range("D39").Value = "'55.59"
$ws.Range("E39").Value = "  +1.64%  "
$ws.Range("D40").Value = "'497.00"
$ws.Range("E40").Value = "  -2.57%  "
$ws.Range("D41").Value = "'0.0429"
$ws.Range("E41").Value = "  +1.66%  "
$ws.Range("D42").Value = "'0.129"
$ws.Range("E42").Value = "  +1.86%  "
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("E45").Value = "  +2.63%  "
$ws.Range("D46").Value = "2.951.23"
$ws.Range("E46").Value = "  +3.65%  "
$ws.Range("D47").Value = "'28.68"
$ws.Range("E47").Value = "  +1.22%  "
$ws.Range("D48").Value = "'2.46"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("E49").Value = "  +2.04%  "
$ws.Range("E51").Value = "  -1.29%  "
